# Insert a new data row at row 359, shifting existing rows 359-429 down to 360-430.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(359).Insert()

$ws.Range("A359").Value = 4
$ws.Range("B359").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C359").Value = "Los Lagos"
$ws.Range("D359").Value = 44995
$ws.Range("E359").Value = 10
$ws.Range("F359").Value = 100112003
$ws.Range("G359").Value = "Ajo"
$ws.Range("H359").Value = "Chilote"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 240
$ws.Range("K359").Value = 20000
$ws.Range("L359").Value = 21000
$ws.Range("M359").Value = 20500
$ws.Range("N359").Value = '$/caja 10 kilos'
$ws.Range("O359").Value = "China"
$ws.Range("P359").Value = 2050
$ws.Range("Q359").Value = 10
$ws.Range("R359").Value = "Hortaliza"
